# Daten aktualisiert am 2024-04-10
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @("MNT-USD", "IMX-USD", "TAO-USD", "GRT-USD", "PEPE-USD")

$startRow = 479
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
